$wb = $excel.ActiveWorkbook

# "Spain" is the closest existing market sheet - use it as the template for
# the new "Turkey" sheet, placing the copy right after it.
$spain = $wb.Worksheets.Item("Spain")
$spain.Copy($null, $spain)

# The freshly created sheet becomes the active sheet, placed right after Spain.
$turkey = $wb.ActiveSheet
$turkey.Name = "Turkey"

# Fill in the Turkey-specific market name / ticket reference.
$turkey.Range("B2").Value = "Turkey Market"
$turkey.Range("B4").Value = "NGC-3191/T3299"

# Column D on the new sheet is a touch narrower than Spain's (25 vs 24.33...).
# The engine's ColumnWidth -> stored width conversion adds a constant offset
# of 5/6 (0.8333...), so back that out to land on exactly 25.
$turkey.Range("D1").EntireColumn.ColumnWidth = 24.1666666666667

# Spain used an enlarged row height (28.8) on rows 3-5 to fit wrapped text;
# Turkey's copy of those rows should fall back to the sheet's default
# (auto-fit) height instead.
$turkey.Range("A3:D5").EntireRow.AutoFit()

# Restore the selection / active-cell bookkeeping: Turkey keeps the cell
# that was selected, Spain's selection becomes the full data range now that
# it is no longer the active tab.
$null = $turkey.Range("G15").Select()
$null = $spain.Range("A1:D12").Select()
$null = $turkey.Select()
